$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-blank D and F columns with new accuracy values,
# copying the percent style used by the adjacent C/E columns.

$ws.Range("D3").Value = 0.9804
$ws.Range("F3").Value = 0.9913

$ws.Range("D4").Value = 0.9677
$ws.Range("F4").Value = 0.9578

$ws.Range("D5").Value = 0.2877
$ws.Range("F5").Value = 0.1443

$ws.Range("D6").Value = 0.1784
$ws.Range("F6").Value = 0.2078

$ws.Range("D7").Value = 0.9657
$ws.Range("F7").Value = 0.992

# Match the percent number format / styling already used on column C and E
$ws.Range("D3:D7").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("F3:F7").NumberFormat = $ws.Range("C3").NumberFormat
$ws.Range("D3:D7").HorizontalAlignment = $ws.Range("C3").HorizontalAlignment
$ws.Range("F3:F7").HorizontalAlignment = $ws.Range("C3").HorizontalAlignment

# Move the active selection to F7, matching the final cell edited
$ws.Range("F7").Select()
